$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2015
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 75
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 38
$ws.Range("N2").Value = -2726

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2506.04
$ws.Range("I131").Value = 2103.0454
$ws.Range("J131").Value = 5461.3335
$ws.Range("K131").Value = 6309.1362
$ws.Range("L131").Value = 16384.0005
$ws.Range("M131").Value = -1269.1362
$ws.Range("N131").Value = -26464.0005

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3847.4827
$ws.Range("I137").Value = 2469.5789
$ws.Range("K137").Value = 7408.736699999999
$ws.Range("M137").Value = -4858.736699999999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9163.65
$ws.Range("I32").Value = 7856.4736
$ws.Range("K32").Value = 7856.4736
$ws.Range("M32").Value = -7569.4736

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 501209.9
$ws.Range("I110").Value = 715269.9
$ws.Range("K110").Value = 715269.9
$ws.Range("M110").Value = -713224.9

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6823.3076
$ws.Range("I132").Value = 6423.273
$ws.Range("J132").Value = 7116.6665
$ws.Range("K132").Value = 19269.819
$ws.Range("L132").Value = 21349.9995
$ws.Range("M132").Value = -16739.819
$ws.Range("N132").Value = -26409.9995

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 105.4
$ws.Range("I22").Value = 31.75
$ws.Range("K22").Value = 31.75
$ws.Range("M22").Value = 141.25

# BSM row 24
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 2192.3635
$ws.Range("I24").Value = 808
$ws.Range("K24").Value = 808
$ws.Range("M24").Value = -573

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 284
$ws.Range("I22").Value = 327.83334
$ws.Range("J22").Value = 152.5
$ws.Range("K22").Value = 327.83334
$ws.Range("L22").Value = 152.5
$ws.Range("M22").Value = 22.16665999999998
$ws.Range("N22").Value = -852.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 419966.5
$ws.Range("I31").Value = 557575.0600000001
$ws.Range("K31").Value = 557575.0600000001
$ws.Range("M31").Value = -557280.0600000001

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 419966.5
$ws.Range("I34").Value = 557575.0600000001
$ws.Range("K34").Value = 557575.0600000001
$ws.Range("M34").Value = -557373.0600000001

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7664.8237
$ws.Range("I99").Value = 8505.833000000001
$ws.Range("K99").Value = 8505.833000000001
$ws.Range("M99").Value = -7007.833000000001

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7664.8237
$ws.Range("I126").Value = 8505.833000000001
$ws.Range("K126").Value = 25517.499
$ws.Range("M126").Value = -23047.499

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3904
$ws.Range("I132").Value = 2648.1072
$ws.Range("K132").Value = 7944.321599999999
$ws.Range("M132").Value = -5414.321599999999

# CUL row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 61149.43
$ws.Range("J36").Value = 106539.25
$ws.Range("L36").Value = 319617.75
$ws.Range("N36").Value = -319955.75

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 500590.94
$ws.Range("I92").Value = 833697.3
$ws.Range("J92").Value = 931.375
$ws.Range("K92").Value = 2501091.9
$ws.Range("L92").Value = 2794.125
$ws.Range("M92").Value = -2499843.9
$ws.Range("N92").Value = -5290.125

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3869.4
$ws.Range("I132").Value = 804
$ws.Range("K132").Value = 7236
$ws.Range("M132").Value = -4706

# GSM row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 27860
$ws.Range("J15").Value = 27860
$ws.Range("L15").Value = 27860
$ws.Range("N15").Value = -28436

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 63332.332
$ws.Range("I24").Value = 70000
$ws.Range("J24").Value = 49997
$ws.Range("K24").Value = 70000
$ws.Range("L24").Value = 49997
$ws.Range("M24").Value = -69827
$ws.Range("N24").Value = -50343

# GSM row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 27860
$ws.Range("J81").Value = 27860
$ws.Range("L81").Value = 27860
$ws.Range("N81").Value = -29856

# GSM row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 27860
$ws.Range("J84").Value = 27860
$ws.Range("L84").Value = 83580
$ws.Range("N84").Value = -93564

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 278538.03
$ws.Range("I132").Value = 360273.7
$ws.Range("K132").Value = 1080821.1
$ws.Range("M132").Value = -1078291.1

# LTW row 4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3549.5
$ws.Range("I4").Value = 2099
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 2099
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -1986
$ws.Range("N4").Value = -5226

# LTW row 12
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2042
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3340

# LTW row 21
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 202.93333
$ws.Range("I21").Value = 181.71428
$ws.Range("J21").Value = 500
$ws.Range("K21").Value = 181.71428
$ws.Range("L21").Value = 500
$ws.Range("M21").Value = -7.714280000000002
$ws.Range("N21").Value = -848

# LTW row 23
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 375166.5
$ws.Range("I23").Value = 487749.75
$ws.Range("K23").Value = 487749.75
$ws.Range("M23").Value = -487519.75

# LTW row 28
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 3549.5
$ws.Range("I28").Value = 2099
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 2099
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -1867
$ws.Range("N28").Value = -5464

# LTW row 31
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 732.2727
$ws.Range("I31").Value = 621.125
$ws.Range("K31").Value = 621.125
$ws.Range("M31").Value = -373.125

# LTW row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 7922.778
$ws.Range("I35").Value = 10500
$ws.Range("J35").Value = 7186.4287
$ws.Range("K35").Value = 10500
$ws.Range("L35").Value = 7186.4287
$ws.Range("M35").Value = -10164
$ws.Range("N35").Value = -7858.4287

# LTW row 37
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 3549.5
$ws.Range("I37").Value = 2099
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 2099
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -1992
$ws.Range("N37").Value = -5214

# LTW row 43
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1264284.9
$ws.Range("J43").Value = 1416665.9
$ws.Range("L43").Value = 1416665.9
$ws.Range("N43").Value = -1417051.9

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4390
$ws.Range("I68").Value = 4150
$ws.Range("J68").Value = 4550
$ws.Range("K68").Value = 4150
$ws.Range("L68").Value = 4550
$ws.Range("M68").Value = -3401
$ws.Range("N68").Value = -6048

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4390
$ws.Range("I71").Value = 4150
$ws.Range("J71").Value = 4550
$ws.Range("K71").Value = 20750
$ws.Range("L71").Value = 22750
$ws.Range("M71").Value = -17006
$ws.Range("N71").Value = -30238

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3591.9473
$ws.Range("I132").Value = 2966.6072
$ws.Range("K132").Value = 8899.821599999999
$ws.Range("M132").Value = -6369.821599999999

# WVR row 7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1000000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1000000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1000000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1000226

# WVR row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 9999
$ws.Range("J28").Value = 9999
$ws.Range("L28").Value = 9999
$ws.Range("N28").Value = -10695

# WVR row 37
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 10029
$ws.Range("J37").Value = 10029
$ws.Range("L37").Value = 10029
